# Update the division problems in the table to the new values from the
# latest generated output (commit 4250d90).
#
# Each pair is applied in the same order the strings occur in the
# document so that a replacement which happens to introduce text that
# equals an *older* (not-yet-processed) search target cannot be matched
# again downstream (e.g. "427÷7=" becomes "650÷7=", but the original
# "650÷7=" earlier in the doc is replaced first).

$d = $word.ActiveDocument

$replacements = @(
    @("420÷4=", "175÷9="),
    @("746÷5=", "480÷3="),
    @("544÷5=", "272÷6="),
    @("650÷7=", "531÷5="),
    @("715÷2=", "928÷2="),
    @("478÷9=", "782÷8="),
    @("349÷6=", "663÷4="),
    @("284÷4=", "417÷5="),
    @("938÷9=", "588÷7="),
    @("895÷7=", "814÷8="),
    @("872÷5=", "797÷9="),
    @("365÷2=", "545÷8="),
    @("535÷7=", "147÷4="),
    @("581÷5=", "227÷8="),
    @("369÷5=", "891÷7="),
    @("489÷6=", "764÷2="),
    @("427÷7=", "650÷7="),
    @("372÷7=", "903÷3="),
    @("168÷7=", "307÷6="),
    @("282÷5=", "680÷7="),
    @("719÷8=", "641÷4="),
    @("734÷3=", "706÷8="),
    @("235÷3=", "212÷3="),
    @("889÷2=", "425÷9="),
    @("875÷4=", "249÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $rng = $d.Content
    $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 1)
}

Write-Host "Replaced $($replacements.Count) division problems."
